$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - update "想去人数" (column F) values for rows 5-12
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 2994
$ws1.Range("F6").Value = 2018
$ws1.Range("F7").Value = 391
$ws1.Range("F8").Value = 138
$ws1.Range("F9").Value = 1117
$ws1.Range("F10").Value = 206
$ws1.Range("F11").Value = 684
$ws1.Range("F12").Value = 64

# Sheet "全部类型" (sheet4) - update "想去人数" (column F) values for rows 5-13
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 2994
$ws4.Range("F6").Value = 2018
$ws4.Range("F7").Value = 391
$ws4.Range("F9").Value = 138
$ws4.Range("F10").Value = 1117
$ws4.Range("F11").Value = 206
$ws4.Range("F12").Value = 684
$ws4.Range("F13").Value = 64
